$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.166052460670471
$ws.Range("B1").Value = 2.428537130355835
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.372776985168457
$ws.Range("E1").Value = 1.23521363735199
